$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20: test.pgm
$ws.Range("A20").Value = "test.pgm"
$ws.Range("B20").Value = 1228800
$ws.Range("C20").Value = 1230848
$ws.Range("D20").Value = 1163232
$ws.Range("E20").Value = 1361920
$ws.Range("F20").Value = 1156096
$ws.Range("G20").Formula = "=AVERAGE(B20:F20)"
$ws.Range("G20").Font.Bold = $true

# New row 21: test_large.pgm
$ws.Range("A21").Value = "test_large.pgm"
$ws.Range("B21").Value = 22693696
$ws.Range("C21").Value = 24498176
$ws.Range("D21").Value = 24958976
$ws.Range("E21").Value = 26193920
$ws.Range("F21").Value = 26352480
$ws.Range("G21").Formula = "=AVERAGE(B21:F21)"
$ws.Range("G21").Font.Bold = $true

# Adjust column A width to fit new longer text (closest reachable value to the
# original author's 14.42578125 stored width, given this engine's internal
# rounding of ColumnWidth to 1/6-character increments)
$ws.Columns.Item(1).ColumnWidth = 13.6666666666667

# Update selection to match diff
$ws.Range("B20:G21").Select()
